$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from the last existing data row (73) down into the
# new row (74) so the new cells inherit the same date / centered-number
# styles, then write the new day's values.
$ws.Range("A73:F73").Copy()
$ws.Range("A74:F74").PasteSpecial(-4122)

$ws.Range("A74").Value = 43976
$ws.Range("B74").Value = 596
$ws.Range("C74").Value = 239
$ws.Range("D74").Value = 400
$ws.Range("E74").Value = 26
$ws.Range("F74").Value = 20

# Grow the "Condicion_Pacientes" table (and its autofilter) to include the
# new row.
$lo = $ws.ListObjects.Item("Condicion_Pacientes")
$lo.Resize($ws.Range("A1:F74"))

# Match the author's saved selection (active cell moved to the new last row).
$null = $ws.Range("F74").Select()
